# Fixed update to excel issue
# - Rename the "Requested quantity" header on the Weekly/Monthly sheets to
#   Weekly_PO_Qty / Monthly_PO_Qty respectively.
# - Add a new "PO Forecast" sheet (ds / PO_Forecast / yhat_lower / yhat_upper)
#   with the forecast rows, styled like the existing sheets.

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new sheet after the last existing sheet so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "PO Forecast"

$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

$newSheet.Cells.Item(2, 1).Value = 45354.99999999999
$newSheet.Cells.Item(2, 2).Value = 15
$newSheet.Cells.Item(2, 3).Value = -27.6938811150257
$newSheet.Cells.Item(2, 4).Value = 60.31841198960385
$newSheet.Cells.Item(3, 1).Value = 45361.99999999999
$newSheet.Cells.Item(3, 2).Value = 16
$newSheet.Cells.Item(3, 3).Value = -30.0781172866651
$newSheet.Cells.Item(3, 4).Value = 56.82863886155312
$newSheet.Cells.Item(4, 1).Value = 45375.99999999999
$newSheet.Cells.Item(4, 2).Value = 18
$newSheet.Cells.Item(4, 3).Value = -29.45663883653351
$newSheet.Cells.Item(4, 4).Value = 62.3848745108851
$newSheet.Cells.Item(5, 1).Value = 45382.99999999999
$newSheet.Cells.Item(5, 2).Value = 19
$newSheet.Cells.Item(5, 3).Value = -26.00810188710268
$newSheet.Cells.Item(5, 4).Value = 67.51034368536649
$newSheet.Cells.Item(6, 1).Value = 45410.99999999999
$newSheet.Cells.Item(6, 2).Value = 23
$newSheet.Cells.Item(6, 3).Value = -22.8696901504698
$newSheet.Cells.Item(6, 4).Value = 68.25498042692965
$newSheet.Cells.Item(7, 1).Value = 45417.99999999999
$newSheet.Cells.Item(7, 2).Value = 24
$newSheet.Cells.Item(7, 3).Value = -21.36334130071876
$newSheet.Cells.Item(7, 4).Value = 68.4826116098278
$newSheet.Cells.Item(8, 1).Value = 45424.99999999999
$newSheet.Cells.Item(8, 2).Value = 25
$newSheet.Cells.Item(8, 3).Value = -22.21809379576028
$newSheet.Cells.Item(8, 4).Value = 71.66977436278627
$newSheet.Cells.Item(9, 1).Value = 45431.99999999999
$newSheet.Cells.Item(9, 2).Value = 26
$newSheet.Cells.Item(9, 3).Value = -21.24893401292222
$newSheet.Cells.Item(9, 4).Value = 72.03481854717646
$newSheet.Cells.Item(10, 1).Value = 45438.99999999999
$newSheet.Cells.Item(10, 2).Value = 27
$newSheet.Cells.Item(10, 3).Value = -17.39153335080346
$newSheet.Cells.Item(10, 4).Value = 75.0148259957776
$newSheet.Cells.Item(11, 1).Value = 45445.99999999999
$newSheet.Cells.Item(11, 2).Value = 28
$newSheet.Cells.Item(11, 3).Value = -16.3493415816246
$newSheet.Cells.Item(11, 4).Value = 73.04991527656392
$newSheet.Cells.Item(12, 1).Value = 45452.99999999999
$newSheet.Cells.Item(12, 2).Value = 29
$newSheet.Cells.Item(12, 3).Value = -16.19825929837553
$newSheet.Cells.Item(12, 4).Value = 73.70012585665802
$newSheet.Cells.Item(13, 1).Value = 45466.99999999999
$newSheet.Cells.Item(13, 2).Value = 31
$newSheet.Cells.Item(13, 3).Value = -12.72227546969722
$newSheet.Cells.Item(13, 4).Value = 74.78626308368395
$newSheet.Cells.Item(14, 1).Value = 45473.99999999999
$newSheet.Cells.Item(14, 2).Value = 32
$newSheet.Cells.Item(14, 3).Value = -14.30535167208299
$newSheet.Cells.Item(14, 4).Value = 75.78387015964312
$newSheet.Cells.Item(15, 1).Value = 45480.99999999999
$newSheet.Cells.Item(15, 2).Value = 33
$newSheet.Cells.Item(15, 3).Value = -13.32490021032511
$newSheet.Cells.Item(15, 4).Value = 81.61088102948506
$newSheet.Cells.Item(16, 1).Value = 45487.99999999999
$newSheet.Cells.Item(16, 2).Value = 34
$newSheet.Cells.Item(16, 3).Value = -8.509927896175476
$newSheet.Cells.Item(16, 4).Value = 78.43870406585219
$newSheet.Cells.Item(17, 1).Value = 45508.99999999999
$newSheet.Cells.Item(17, 2).Value = 37
$newSheet.Cells.Item(17, 3).Value = -10.30796772943678
$newSheet.Cells.Item(17, 4).Value = 84.29940904321295
$newSheet.Cells.Item(18, 1).Value = 45515.99999999999
$newSheet.Cells.Item(18, 2).Value = 38
$newSheet.Cells.Item(18, 3).Value = -6.513548388273168
$newSheet.Cells.Item(18, 4).Value = 84.50042325483415
$newSheet.Cells.Item(19, 1).Value = 45529.99999999999
$newSheet.Cells.Item(19, 2).Value = 40
$newSheet.Cells.Item(19, 3).Value = -1.761048943074307
$newSheet.Cells.Item(19, 4).Value = 87.60553282083502
$newSheet.Cells.Item(20, 1).Value = 45571.99999999999
$newSheet.Cells.Item(20, 2).Value = 46
$newSheet.Cells.Item(20, 3).Value = -0.1717543417316129
$newSheet.Cells.Item(20, 4).Value = 94.83261955306494
$newSheet.Cells.Item(21, 1).Value = 45578.99999999999
$newSheet.Cells.Item(21, 2).Value = 47
$newSheet.Cells.Item(21, 3).Value = 3.955307069805586
$newSheet.Cells.Item(21, 4).Value = 91.07953516018536
$newSheet.Cells.Item(22, 1).Value = 45585.99999999999
$newSheet.Cells.Item(22, 2).Value = 48
$newSheet.Cells.Item(22, 3).Value = 2.142861708537002
$newSheet.Cells.Item(22, 4).Value = 93.80507190276657
$newSheet.Cells.Item(23, 1).Value = 45592.99999999999
$newSheet.Cells.Item(23, 2).Value = 49
$newSheet.Cells.Item(23, 3).Value = 2.296477742153588
$newSheet.Cells.Item(23, 4).Value = 91.25445083878246
$newSheet.Cells.Item(24, 1).Value = 45599.99999999999
$newSheet.Cells.Item(24, 2).Value = 50
$newSheet.Cells.Item(24, 3).Value = 5.39299268440499
$newSheet.Cells.Item(24, 4).Value = 99.12342070886315
$newSheet.Cells.Item(25, 1).Value = 45606.99999999999
$newSheet.Cells.Item(25, 2).Value = 51
$newSheet.Cells.Item(25, 3).Value = 6.380352850412285
$newSheet.Cells.Item(25, 4).Value = 96.87049287287786
$newSheet.Cells.Item(26, 1).Value = 45613.99999999999
$newSheet.Cells.Item(26, 2).Value = 52
$newSheet.Cells.Item(26, 3).Value = 5.00343181910051
$newSheet.Cells.Item(26, 4).Value = 98.82579890061714
$newSheet.Cells.Item(27, 1).Value = 45620.99999999999
$newSheet.Cells.Item(27, 2).Value = 53
$newSheet.Cells.Item(27, 3).Value = 6.520188673345296
$newSheet.Cells.Item(27, 4).Value = 97.36116270027058
$newSheet.Cells.Item(28, 1).Value = 45627.99999999999
$newSheet.Cells.Item(28, 2).Value = 55
$newSheet.Cells.Item(28, 3).Value = 8.664495880557915
$newSheet.Cells.Item(28, 4).Value = 103.5370778682293
$newSheet.Cells.Item(29, 1).Value = 45634.99999999999
$newSheet.Cells.Item(29, 2).Value = 56
$newSheet.Cells.Item(29, 3).Value = 8.572874615925379
$newSheet.Cells.Item(29, 4).Value = 101.1188301324053
$newSheet.Cells.Item(30, 1).Value = 45641.99999999999
$newSheet.Cells.Item(30, 2).Value = 57
$newSheet.Cells.Item(30, 3).Value = 12.40345009506207
$newSheet.Cells.Item(30, 4).Value = 99.49733173483644

# Match the formatting used on the existing sheets: bold/centered/bordered
# header row, and the date/time number format on column A's data rows.
$wsWeekly.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$newSheet.Range("A2:A30").PasteSpecial(-4122)

$excel.CutCopyMode = 0
